$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "7級"
$wb.Worksheets.Item(2).Name = "7級用対策プリント作成シート"
$wb.Worksheets.Item(3).Name = "6級"
$wb.Worksheets.Item(4).Name = "6級用対策プリント作成シート"
$wb.Worksheets.Item(5).Name = "5級"
$wb.Worksheets.Item(6).Name = "5級用対策プリント作成シート"
